$p = $ppt.ActivePresentation

# --- 1) Bump the cached "datetimeFigureOut" footer date from 7/30/17 to
#        5/25/18 everywhere it is rendered: every slide layout, the slide
#        master, and the notes master each carry their own "Date
#        Placeholder" shape with the old cached text. ---

$oldDate = "7/30/17"
$newDate = "5/25/18"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# --- 2) Fix the typo in puzzle 3 (slide 5): the expected output was
#        "1 2 -1" but should read "1 6 -1". ---

$slide = $p.Slides.Item(5)
$puzzleShape = $slide.Shapes.Item(2)
$puzzleRange = $puzzleShape.TextFrame.TextRange
$target = $puzzleRange.Find("2 ")
if ($target -ne $null) {
    $target.Text = "6 "
}
